# "added further analysis points"
#
# Adds two new trailing columns to the export:
#   F  "Årsag"          (reason the customer cancelled) - replaces the old
#                        "TCV_range" column, which is pushed out to column H
#   G  "Ny leverandør"  (new supplier, only populated for one row)
#   H  "TCV_range"      (the bucket value that used to live in column F)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- header row -----------------------------------------------------------
$ws.Range("F1").Value = "Årsag"
$ws.Range("G1").Value = "Ny leverandør"
$ws.Range("H1").Value = "TCV_range"

# Give the two new header cells the same look (bold, centered, bordered)
# as the rest of the header row.
$ws.Range("F1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)   # xlPasteFormats

# --- data rows -------------------------------------------------------------
$tcvRange = "140000-160000"

$reasons = @{
    2  = "Fusionerer med anden virksomhed"
    3  = "Utilfredshed (Service - uddyb i bemærkninger)"
    4  = "Fusionerer med anden virksomhed"
    5  = "Utilfredshed (Service - uddyb i bemærkninger)"
    6  = "Ikke oplyst"
    7  = "Utilfredshed (Service - uddyb i bemærkninger)"
    8  = "Konkurs"
    9  = "Ikke oplyst"
    10 = "Bruger ikke produktet"
    11 = "Strategisk beslutning"
    12 = "Strategisk beslutning"
    13 = "Ikke oplyst"
    14 = "Ikke oplyst"
    15 = "Insourcing af lønnen (anden leverandør)"
}

foreach ($row in 2..15) {
    $ws.Range("F$row").Value = $reasons[$row]
    $ws.Range("H$row").Value = $tcvRange
}

# Only one row records a new supplier.
$ws.Range("G13").Value = "DataLøn"
